$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds month-start dates used as column headers. Several of them
# (F1, G1, H1, I1, Z1, AB1, AC1) had been entered/left as plain text
# (" Set/23", " Out/23", " Nov/23", " Dez/23", " Mai/25", " Jul/25", " Ago/25")
# instead of real date values, unlike their neighbouring cells. Fix them by
# writing the correct date serial values so Excel stores them as dates
# (using the same [$-416]mmm-yy;@ custom format already applied to the row).
$ws.Range("F1").Value = 45170   # 01/09/2023
$ws.Range("G1").Value = 45200   # 01/10/2023
$ws.Range("H1").Value = 45231   # 01/11/2023
$ws.Range("I1").Value = 45261   # 01/12/2023
$ws.Range("Z1").Value = 45778   # 01/05/2025
$ws.Range("AB1").Value = 45839  # 01/07/2025
$ws.Range("AC1").Value = 45870  # 01/08/2025

# A1 (top-left corner cell above the row labels) picks up the same
# date number format with left alignment used elsewhere in the header row.
$ws.Range("A1").NumberFormat = "[$-416]mmm\-yy;@"
$ws.Range("A1").HorizontalAlignment = -4131

# Update the window/view state: scroll the visible area to the right and
# make Z2 the active/selected cell.
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Z2").Select()

Write-Host "done"
